$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.0003978386713523962
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2.016951252558324
$ws.Range("H3").Value = 0.2031221295889138

$ws.Range("C4").Value = 0.0003978386713523962
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

$ws.Range("C6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

$ws.Range("C7").Value = 2.016951252558324
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

$ws.Range("C8").Value = 0.2031221295889138
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("I8").Value = 0

$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 0
